$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the "Temps réel (j)" (column D) values for the task rows
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 2.75
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 2.75
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 0

# Move active selection to D3, matching the workbook's saved cursor position
$ws.Range("D3").Select()
